$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("D3").Value = "AUTOREJECT"
$ws.Range("E3").Value = "Threshold"

# Column D: AUTOREJECT fraction-of-epochs-dropped values (percent formatted)
$ws.Range("D4").Value = [double]"0.8"
$ws.Range("D5").Value = [double]"0"
$ws.Range("D6").Value = [double]"0"
$ws.Range("D7").Value = [double]"0"
$ws.Range("D8").Value = [double]"0"
$ws.Range("D9").Value = [double]"0"
$ws.Range("D10").Value = [double]"0"
$ws.Range("D11").Value = [double]"0"
$ws.Range("D12").Value = [double]"5.0000000000000001E-3"
$ws.Range("D13").Value = [double]"0"
$ws.Range("D14").Value = [double]"0"
$ws.Range("D15").Value = [double]"0"
$ws.Range("D16").Value = [double]"0"
$ws.Range("D17").Value = [double]"0.2"
$ws.Range("D18").Value = [double]"0.4"
$ws.Range("D19").Value = [double]"0"
$ws.Range("D20").Value = [double]"0.6"
$ws.Range("D21").Value = [double]"0.4"
$ws.Range("D22").Value = [double]"0.2"
$ws.Range("D23").Value = [double]"0.4"
$ws.Range("D24").Value = [double]"4.5999999999999999E-2"
$ws.Range("D25").Value = [double]"1E-3"
$ws.Range("D26").Value = [double]"0"
$ws.Range("D27").Value = [double]"0"
$ws.Range("D28").Value = [double]"0.6"
$ws.Range("D29").Value = [double]"0"
$ws.Range("D30").Value = [double]"0"
$ws.Range("D31").Value = [double]"0"
$ws.Range("D32").Value = [double]"0"
$ws.Range("D33").Value = [double]"0"
$ws.Range("D34").Value = [double]"0"
$ws.Range("D35").Value = [double]"0.33299999999999996"
$ws.Range("D36").Value = [double]"0.2"
$ws.Range("D37").Value = [double]"0"
$ws.Range("D38").Value = [double]"0.2"
$ws.Range("D39").Value = [double]"0"
$ws.Range("D40").Value = [double]"1E-3"
$ws.Range("D41").Value = [double]"0"
$ws.Range("D42").Value = [double]"0"
$ws.Range("D43").Value = [double]"0.6"
$ws.Range("D44").Value = [double]"0.2"
$ws.Range("D45").Value = [double]"0"
$ws.Range("D46").Value = [double]"0.312"
$ws.Range("D47").Value = [double]"0.8"
$ws.Range("D48").Value = [double]"0"
$ws.Range("D49").Value = [double]"0"
$ws.Range("D50").Value = [double]"0"
$ws.Range("D51").Value = [double]"0"
$ws.Range("D52").Value = [double]"0.8"
$ws.Range("D53").Value = [double]"0.6"

# Column E: Threshold values (plain numeric, a few in scientific notation)
$ws.Range("E4").Value = [double]"3.2868127119523401E-5"
$ws.Range("E5").Value = [double]"1.1354786009950301E-3"
$ws.Range("E6").Value = [double]"3.3041984184570801E-3"
$ws.Range("E7").Value = [double]"2.05740405048585E-4"
$ws.Range("E8").Value = [double]"1.4738356454772701E-3"
$ws.Range("E9").Value = [double]"1.4314424280611701E-3"
$ws.Range("E10").Value = [double]"5.91921310721786E-4"
$ws.Range("E11").Value = [double]"3.1069460573752601E-3"
$ws.Range("E12").Value = [double]"5.1330900475964904E-4"
$ws.Range("E13").Value = [double]"5.4487834129254701E-4"
$ws.Range("E14").Value = [double]"1.9002787618916501E-3"
$ws.Range("E15").Value = [double]"4.6730206729784698E-3"
$ws.Range("E16").Value = [double]"0.18602805874174599"
$ws.Range("E17").Value = [double]"5.0019790055730198E-5"
$ws.Range("E18").Value = [double]"2.3418901681103101E-3"
$ws.Range("E19").Value = [double]"2.40286305161034E-3"
$ws.Range("E20").Value = [double]"1.5279868748495401E-3"
$ws.Range("E21").Value = [double]"1.0804329792613501E-3"
$ws.Range("E22").Value = [double]"2.8188762723774302E-3"
$ws.Range("E23").Value = [double]"7.1360624853334905E-4"
$ws.Range("E24").Value = [double]"1.27985369720505E-4"
$ws.Range("E25").Value = [double]"5.4443869150731805E-4"
$ws.Range("E26").Value = [double]"1.0509863830914901E-2"
$ws.Range("E27").Value = [double]"2.70212461224029E-4"
$ws.Range("E28").Value = [double]"9.21856458575309E-4"
$ws.Range("E29").Value = [double]"7.76994632038501E-4"
$ws.Range("E30").Value = [double]"2.9474303542584099E-2"
$ws.Range("E31").Value = [double]"4.5137187057741398E-4"
$ws.Range("E32").Value = [double]"9.2237056727030899E-3"
$ws.Range("E33").Value = [double]"1.6538274595367501E-3"
$ws.Range("E43").Value = [double]"2.9658319767110898E-5"
$ws.Range("E44").Value = [double]"2.1335968638865101E-4"
$ws.Range("E45").Value = [double]"3.2944759755447199E-3"
$ws.Range("E46").Value = [double]"2.6728403358308402E-4"
$ws.Range("E47").Value = [double]"3.7365993376928498E-5"
$ws.Range("E48").Value = [double]"3.3466163271822399E-4"
$ws.Range("E49").Value = [double]"6.4712009797319298E-4"
$ws.Range("E50").Value = [double]"4.8302818035576899E-4"
$ws.Range("E51").Value = [double]"1.80003919940958E-2"
$ws.Range("E52").Value = [double]"1.1170761088547801E-3"
$ws.Range("E53").Value = [double]"6.5361291315035002E-5"

# Apply percent number format (0.00%) to the whole AUTOREJECT column of data
$ws.Range("D4:D53").NumberFormat = "0.00%"

# A handful of very small Threshold values are shown in scientific notation
foreach ($addr in @("E17","E43","E44","E45","E46","E47","E53")) {
    $ws.Range($addr).NumberFormat = "0.00E+00"
}

# Stray formatted (empty) cell that appears in the edited range
$ws.Range("H5").HorizontalAlignment = -4108

# Column widths for the two new columns
$ws.Columns.Item(4).ColumnWidth = 10.8
$ws.Columns.Item(5).ColumnWidth = 22.6

# Restore the selection/scroll state used while authoring the sheet
$ws.Range("G48").Select()
$excel.ActiveWindow.ScrollRow = 25
